# Apply weekly price updates to "Fruta, Macroferia Regional de Talca - Mandarina"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 163: Clemenuless/Primera -> Murcott/Especial, date + price refresh ---
$ws.Range("D163").Value = 44448
$ws.Range("K163").Value = "Murcott"
$ws.Range("L163").Value = "Especial"
$ws.Range("M163").Value = 200
$ws.Range("N163").Value = 8000
$ws.Range("O163").Value = 8000
$ws.Range("P163").Value = 8000
$ws.Range("S163").Value = 800

# --- Row 164: Clemenuless -> Murcott, date refresh ---
$ws.Range("D164").Value = 44448
$ws.Range("K164").Value = "Murcott"

# --- Row 165: date + price refresh, unit changes to bandeja ---
$ws.Range("D165").Value = 44399
$ws.Range("M165").Value = 240
$ws.Range("N165").Value = 6000
$ws.Range("O165").Value = 6000
$ws.Range("P165").Value = 6000
$ws.Range("Q165").Value = "$/bandeja 10 kilos"
$ws.Range("S165").Value = 600

# --- Row 166: Segunda -> Primera, price refresh, unit changes to bandeja ---
$ws.Range("L166").Value = "Primera"
$ws.Range("M166").Value = 250
$ws.Range("N166").Value = 6000
$ws.Range("O166").Value = 6000
$ws.Range("P166").Value = 6000
$ws.Range("Q166").Value = "$/bandeja 10 kilos"
$ws.Range("S166").Value = 600

# --- New row 167: Clemenuless / Primera ---
$ws.Range("A167").Value = 5
$ws.Range("B167").Value = "Macroferia Regional de Talca"
$ws.Range("C167").Value = "Maule"
$ws.Range("D167").Value = 44400
$ws.Range("D167").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E167").Value = 7
$ws.Range("F167").Value = "Fruta"
$ws.Range("G167").Value = 100102
$ws.Range("H167").Value = "Cítricos"
$ws.Range("I167").Value = 100102004
$ws.Range("J167").Value = "Mandarina"
$ws.Range("K167").Value = "Clemenuless"
$ws.Range("L167").Value = "Primera"
$ws.Range("M167").Value = 230
$ws.Range("N167").Value = 5000
$ws.Range("O167").Value = 5000
$ws.Range("P167").Value = 5000
$ws.Range("Q167").Value = "$/caja 10 kilos"
$ws.Range("R167").Value = "Provincia de Limarí"
$ws.Range("S167").Value = 500
$ws.Range("T167").Value = 10

# --- New row 168: Clemenuless / Segunda ---
$ws.Range("A168").Value = 5
$ws.Range("B168").Value = "Macroferia Regional de Talca"
$ws.Range("C168").Value = "Maule"
$ws.Range("D168").Value = 44400
$ws.Range("D168").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E168").Value = 7
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100102
$ws.Range("H168").Value = "Cítricos"
$ws.Range("I168").Value = 100102004
$ws.Range("J168").Value = "Mandarina"
$ws.Range("K168").Value = "Clemenuless"
$ws.Range("L168").Value = "Segunda"
$ws.Range("M168").Value = 200
$ws.Range("N168").Value = 4000
$ws.Range("O168").Value = 4000
$ws.Range("P168").Value = 4000
$ws.Range("Q168").Value = "$/caja 10 kilos"
$ws.Range("R168").Value = "Provincia de Limarí"
$ws.Range("S168").Value = 400
$ws.Range("T168").Value = 10

# --- New row 169: Clemenuless / Tercera (shifted down from the old row 167) ---
$ws.Range("A169").Value = 5
$ws.Range("B169").Value = "Macroferia Regional de Talca"
$ws.Range("C169").Value = "Maule"
$ws.Range("D169").Value = 44400
$ws.Range("D169").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100102
$ws.Range("H169").Value = "Cítricos"
$ws.Range("I169").Value = 100102004
$ws.Range("J169").Value = "Mandarina"
$ws.Range("K169").Value = "Clemenuless"
$ws.Range("L169").Value = "Tercera"
$ws.Range("M169").Value = 1003
$ws.Range("N169").Value = 3000
$ws.Range("O169").Value = 3000
$ws.Range("P169").Value = 3000
$ws.Range("Q169").Value = "$/caja 10 kilos"
$ws.Range("R169").Value = "Provincia de Limarí"
$ws.Range("S169").Value = 300
$ws.Range("T169").Value = 10
